$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.710.40"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.162.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.36%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.58"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.09%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.159.90"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.49%  "
$ws.Range("E9").Value = "  -6.61%  "
$ws.Range("E10").Value = "  -8.89%  "
$ws.Range("E11").Value = "  -8.03%  "
$ws.Range("E12").Value = "  -4.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.704.58"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.161.52"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.657.85"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("E18").Value = "  -8.54%  "
$ws.Range("E19").Value = "  -6.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.24%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.515"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.294.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0969"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -10.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("E35").Value = "  -5.81%  "
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.27"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.87%  "
$ws.Range("E39").Value = "  -8.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.41"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0702"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.190.98"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.699"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.96%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.47"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.275.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.70"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.41%  "
